$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Invited Lecturer" block (rows 22-27): re-sort by year, ascending
$ws.Range("A22").Value = "Health Sciences Day, Erasmus MC, Netherlands"
$ws.Range("B22").Value = 2019

$ws.Range("A23").Value = "Leibniz Institute for Prevention Researchand Epidemiology - BIPS, Germany"
$ws.Range("B23").Value = 2019

$ws.Range("A24").Value = "Clinical Research Department, Hospital Italiano de Buenos Aires, Argentina"
$ws.Range("B24").Value = 2020

$ws.Range("A25").Value = "Open Science Community Rotterdam, Netherlands"
$ws.Range("B25").Value = 2020

$ws.Range("A26").Value = "Epidemiology Department, Leiden University Medical Center, Netherlands"
$ws.Range("B26").Value = 2020

$ws.Range("A27").Value = "Tutorium " + [char]0x201C + "An Introduction to Causal Inference and Target Trials" + [char]0x201D + " - IBS"
$ws.Range("B27").Value = 2021

# "Research Visits" block (rows 33-34): re-sort by year, ascending
$ws.Range("A33").Value = "Harvard T.H. Chan School of Public Health, Boston, USA"
$ws.Range("B33").Value = 2018

$ws.Range("A34").Value = "Leibniz Institute for Prevention Researchand Epidemiology, Bremen, Germany"
$ws.Range("B34").Value = 2019

# "Institutional and Community Service" block (rows 36 & 38): re-sort by year, ascending
$ws.Range("A36").Value = "Organizer R-Ladies Rotterdam"
$ws.Range("B36").Value = "2018 - 2020"
$ws.Range("C36").ClearContents()

$ws.Range("A38").Value = "Organizer Epidemiology Seminars "
$ws.Range("B38").Value = "2019 - 2021"
$ws.Range("C38").Value = 2

# Make column B a bit wider so the years fit comfortably
$ws.Columns.Item(2).ColumnWidth = 11.6

# Update the view: zoom in and leave the newly-sorted "Research Visits" /
# "Institutional and Community Service" rows selected, as the user did
# after finishing the re-sort.
$excel.ActiveWindow.Zoom = 157
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A28:D33").Select() | Out-Null
